# Advance the "today" reference date for this tracking sheet from
# 2025-11-04 to 2025-11-05, updating the "剩余" (remaining days) column (E)
# and, where the countdown rolled over (remaining was 1 -> would hit 0),
# resetting the countdown: E becomes the total-days value (D) and the
# "开始时间" (start date, column F) becomes the new reference date
# (20251105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldToday = Get-Date -Year 2025 -Month 11 -Day 4
$newTodayStr = "20251105"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($row = 2; $row -le $lastRow; $row++) {
    $dCell = $ws.Cells.Item($row, 4)   # D: 总天 (total days)
    $eCell = $ws.Cells.Item($row, 5)   # E: 剩余 (remaining days)
    $fCell = $ws.Cells.Item($row, 6)   # F: 开始时间 (start date, yyyymmdd)

    $dVal = $dCell.Value()
    $eVal = $eCell.Value()
    $fVal = $fCell.Value()

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fStr = [string]([int]$fVal)
    if ($fStr.Length -ne 8) {
        # Malformed start date (e.g. "202510929") -- leave this row untouched.
        continue
    }

    $fYear = [int]$fStr.Substring(0, 4)
    $fMonth = [int]$fStr.Substring(4, 2)
    $fDay = [int]$fStr.Substring(6, 2)
    $fDate = Get-Date -Year $fYear -Month $fMonth -Day $fDay

    $elapsed = [int]$dVal - [int]$eVal
    $computedToday = $fDate.AddDays($elapsed)

    if ($computedToday.Date -ne $oldToday.Date) {
        # Row's countdown isn't anchored on the previous reference date;
        # leave it alone.
        continue
    }

    $newE = [int]$eVal - 1
    if ($newE -le 0) {
        $eCell.Value = [int]$dVal
        $fCell.Value = [int]$newTodayStr
    } else {
        $eCell.Value = $newE
    }
}
